# Update cryptos list prices / 1h volume percentages.
# Rows 45 and 46 swap: NEARProtocol moves up to row 45, EnergySwap moves to row 46,
# each carrying its own refreshed price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.587.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.007.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.49"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5026"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4257"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.77"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.128"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.111"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.545"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.926.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -11.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001126"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06660"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.16%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.992"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.606.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.98"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.79"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.561"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.343"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.43"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.82%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -8.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.596"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -9.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09971"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.875"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.787"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.654"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02478"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.312"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06383"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6573"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.79"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2079"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6356"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.27%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.217"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.48%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.44"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.303"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.532"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07031"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.142"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.09%  "
